$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new Streamlit submission as row 24 (dima roman tarabzouni).

$ws.Range("A24").Value = 'dima roman tarabzouni_20251202_133740'
$ws.Range("B24").Value = "'"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = 'dima roman tarabzouni'
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 'Female'
$ws.Range("F24").Value = '2025-12-02 13:37:40'
$ws.Range("G24").Value = "{`n  ""portion"": 0.8,`n  ""diet"": 1.0,`n  ""salt"": 0.4,`n  ""fat"": 0.6,`n  ""natural"": 0.4,`n  ""convenience"": 0.0,`n  ""price"": 1.0`n}"
$ws.Range("H24").Value = 'Nongshim Neoguri Spicy Seafood'
$ws.Range("I24").Value = "'0.576"
$ws.Range("I24").Style = "Normal"
$ws.Range("J24").Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range("K24").Value = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range("L24").Value = "'0.558"
$ws.Range("L24").Style = "Normal"
$ws.Range("M24").Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range("N24").Value = 'Nongshim Shin Ramyun'
$ws.Range("O24").Value = "'0.498"
$ws.Range("O24").Style = "Normal"
$ws.Range("P24").Value = 'Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio'
$ws.Range("Q24").Value = 'Kraft Macaroni & Cheese Dinner'
$ws.Range("R24").Value = "'0.647"
$ws.Range("R24").Style = "Normal"
$ws.Range("S24").Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range("T24").Value = 'Amy’s Macaroni & Cheese (frozen)'
$ws.Range("U24").Value = "'0.570"
$ws.Range("U24").Style = "Normal"
$ws.Range("V24").Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range("W24").Value = 'Annie’s Shells & White Cheddar'
$ws.Range("X24").Value = "'0.567"
$ws.Range("X24").Style = "Normal"
$ws.Range("Y24").Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'
$ws.Range("Z24").Value = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range("AA24").Value = "'0.610"
$ws.Range("AA24").Style = "Normal"
$ws.Range("AB24").Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range("AC24").Value = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range("AD24").Value = "'0.532"
$ws.Range("AD24").Style = "Normal"
$ws.Range("AE24").Value = 'Portátil, saludable, fácil, buena textura, sabor suave'
$ws.Range("AF24").Value = 'Jack Link’s Beef Jerky Original'
$ws.Range("AG24").Value = "'0.501"
$ws.Range("AG24").Style = "Normal"
$ws.Range("AH24").Value = 'Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña'
